# Scheduled-runner refresh: update cached market-price / profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9262397
$ws.Range("I64").Value = 11114216
$ws.Range("J64").Value = 3300
$ws.Range("K64").Value = 11114216
$ws.Range("L64").Value = 3300
$ws.Range("M64").Value = -11113968
$ws.Range("N64").Value = -3796
$ws.Range("H67").Value = 9262397
$ws.Range("I67").Value = 11114216
$ws.Range("J67").Value = 3300
$ws.Range("K67").Value = 11114216
$ws.Range("L67").Value = 3300
$ws.Range("M67").Value = -11113358
$ws.Range("N67").Value = -5016
$ws.Range("H76").Value = 6059.375
$ws.Range("I76").Value = 4706.6665
$ws.Range("J76").Value = 7252.9414
$ws.Range("K76").Value = 4706.6665
$ws.Range("L76").Value = 7252.9414
$ws.Range("M76").Value = -4391.6665
$ws.Range("N76").Value = -7882.9414
$ws.Range("H79").Value = 6059.375
$ws.Range("I79").Value = 4706.6665
$ws.Range("J79").Value = 7252.9414
$ws.Range("K79").Value = 4706.6665
$ws.Range("L79").Value = 7252.9414
$ws.Range("M79").Value = -3614.6665
$ws.Range("N79").Value = -9436.9414
$ws.Range("H82").Value = 1473.3334
$ws.Range("I82").Value = 1473.3334
$ws.Range("K82").Value = 4420.0002
$ws.Range("M82").Value = -4014.0002
$ws.Range("H85").Value = 1473.3334
$ws.Range("I85").Value = 1473.3334
$ws.Range("K85").Value = 4420.0002
$ws.Range("M85").Value = -3016.0002
$ws.Range("H138").Value = 11906482
$ws.Range("I138").Value = 1701.6471
$ws.Range("J138").Value = 62501800
$ws.Range("K138").Value = 5104.9413
$ws.Range("L138").Value = 187505400
$ws.Range("M138").Value = 35.0587000000005
$ws.Range("N138").Value = -187515680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7392.7256
$ws.Range("I32").Value = 7894.096
$ws.Range("K32").Value = 7894.096
$ws.Range("M32").Value = -7607.096
$ws.Range("H132").Value = 8066668
$ws.Range("I132").Value = 10001764
$ws.Range("J132").Value = 3767.8333
$ws.Range("K132").Value = 30005292
$ws.Range("L132").Value = 11303.4999
$ws.Range("M132").Value = -30002762
$ws.Range("N132").Value = -16363.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1116.55
$ws.Range("I20").Value = 1188.7693
$ws.Range("K20").Value = 1188.7693
$ws.Range("M20").Value = -941.7692999999999
$ws.Range("H26").Value = 26575.428
$ws.Range("I26").Value = 8617.75
$ws.Range("K26").Value = 8617.75
$ws.Range("M26").Value = -8325.75
$ws.Range("H28").Value = 30000
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30588
$ws.Range("H87").Value = 38354
$ws.Range("J87").Value = 38354
$ws.Range("L87").Value = 38354
$ws.Range("N87").Value = -40850
$ws.Range("H90").Value = 38354
$ws.Range("J90").Value = 38354
$ws.Range("L90").Value = 115062
$ws.Range("N90").Value = -127542
$ws.Range("H107").Value = 3175.6667
$ws.Range("I107").Value = 3708
$ws.Range("J107").Value = 2111
$ws.Range("K107").Value = 3708
$ws.Range("L107").Value = 2111
$ws.Range("M107").Value = -1788
$ws.Range("N107").Value = -5951
$ws.Range("H134").Value = 6282.6387
$ws.Range("I134").Value = 5626.724
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 16880.172
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -14345.172
$ws.Range("N134").Value = -32070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2293.5293
$ws.Range("J62").Value = 2290
$ws.Range("L62").Value = 2290
$ws.Range("N62").Value = -3538
$ws.Range("H65").Value = 2293.5293
$ws.Range("J65").Value = 2290
$ws.Range("L65").Value = 11450
$ws.Range("N65").Value = -17690
$ws.Range("H107").Value = 567.75
$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 506.14285
$ws.Range("K107").Value = 999
$ws.Range("L107").Value = 506.14285
$ws.Range("M107").Value = 921
$ws.Range("N107").Value = -4346.14285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 208.3
$ws.Range("J38").Value = 67.14286
$ws.Range("L38").Value = 201.42858
$ws.Range("N38").Value = -895.42858
$ws.Range("H56").Value = 3978.261
$ws.Range("I56").Value = 3978.261
$ws.Range("K56").Value = 3978.261
$ws.Range("M56").Value = -3448.261
$ws.Range("H117").Value = 1593.4
$ws.Range("J117").Value = 2163.3333
$ws.Range("L117").Value = 6489.999899999999
$ws.Range("N117").Value = -13373.9999
$ws.Range("H129").Value = 4082.8823
$ws.Range("J129").Value = 5573.5454
$ws.Range("L129").Value = 16720.6362
$ws.Range("N129").Value = -26720.6362
$ws.Range("H134").Value = 3845.4546
$ws.Range("I134").Value = 1808.3334
$ws.Range("J134").Value = 6290
$ws.Range("K134").Value = 5425.0002
$ws.Range("L134").Value = 18870
$ws.Range("M134").Value = -355.0002000000004
$ws.Range("N134").Value = -29010
$ws.Range("H139").Value = 2076.8147
$ws.Range("I139").Value = 1593.8572
$ws.Range("J139").Value = 2596.923
$ws.Range("K139").Value = 4781.571599999999
$ws.Range("L139").Value = 7790.768999999999
$ws.Range("M139").Value = 358.4284000000007
$ws.Range("N139").Value = -18070.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14196757
$ws.Range("I80").Value = 19610428
$ws.Range("K80").Value = 19610428
$ws.Range("M80").Value = -19609430
$ws.Range("H83").Value = 14196757
$ws.Range("I83").Value = 19610428
$ws.Range("K83").Value = 98052140
$ws.Range("M83").Value = -98047148
$ws.Range("H122").Value = 2900163.2
$ws.Range("I122").Value = 4445809.5
$ws.Range("J122").Value = 2076.5
$ws.Range("K122").Value = 13337428.5
$ws.Range("L122").Value = 6229.5
$ws.Range("M122").Value = -13334978.5
$ws.Range("N122").Value = -11129.5
$ws.Range("H132").Value = 4141.793
$ws.Range("I132").Value = 3166.125
$ws.Range("J132").Value = 5342.615
$ws.Range("K132").Value = 9498.375
$ws.Range("L132").Value = 16027.845
$ws.Range("M132").Value = -6968.375
$ws.Range("N132").Value = -21087.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9809.091
$ws.Range("J7").Value = 6628.5713
$ws.Range("L7").Value = 6628.5713
$ws.Range("N7").Value = -6852.5713
$ws.Range("H22").Value = 786.96295
$ws.Range("I22").Value = 748.82355
$ws.Range("J22").Value = 851.8
$ws.Range("K22").Value = 748.82355
$ws.Range("L22").Value = 851.8
$ws.Range("M22").Value = -453.82355
$ws.Range("N22").Value = -1441.8
$ws.Range("H27").Value = 786.96295
$ws.Range("I27").Value = 748.82355
$ws.Range("J27").Value = 851.8
$ws.Range("K27").Value = 748.82355
$ws.Range("L27").Value = 851.8
$ws.Range("M27").Value = -641.82355
$ws.Range("N27").Value = -1065.8
$ws.Range("H46").Value = 950.7646999999999
$ws.Range("I46").Value = 809.1111
$ws.Range("J46").Value = 1110.125
$ws.Range("K46").Value = 809.1111
$ws.Range("L46").Value = 1110.125
$ws.Range("M46").Value = -621.1111
$ws.Range("N46").Value = -1486.125
$ws.Range("H126").Value = 9809.091
$ws.Range("J126").Value = 6628.5713
$ws.Range("L126").Value = 19885.7139
$ws.Range("N126").Value = -24825.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3813.15
$ws.Range("J126").Value = 10400
$ws.Range("L126").Value = 31200
$ws.Range("N126").Value = -36140.001
